$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 9 = ASTUDILLO ESPINOZA JOSE MANUEL: record September sale broken down by product group
$wsGrupo.Range("I9").Value = 572.4     # LAVABOS
$wsGrupo.Range("L9").Value = 537.34    # PIEDRA SINTERIZADA
$wsGrupo.Range("M9").Value = 1814.4    # PORCELANATO

# Row 55 = count of clients with sales "N de 53" per product group, bump by 1 for the columns touched
$wsGrupo.Range("I55").Value = "7 de 53"
$wsGrupo.Range("L55").Value = "5 de 53"
$wsGrupo.Range("M55").Value = "11 de 53"

# --- Sheet 2: "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 9 = ASTUDILLO ESPINOZA JOSE MANUEL, column F = septiembre
$wsMensual.Range("F9").Value = 2924.14

# Row 59 = totals row
$wsMensual.Range("F59").Value = 26565.33

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 7 = LAVABOS
$wsCumplimiento.Range("D7").Value = 1622.7
$wsCumplimiento.Range("E7").Value = -735.988983712426
$wsCumplimiento.Range("F7").Value = 1.83002124727605

# Row 11 = PIEDRA SINTERIZADA
$wsCumplimiento.Range("D11").Value = 6612.96
$wsCumplimiento.Range("E11").Value = 11218.4543984654
$wsCumplimiento.Range("F11").Value = 0.370860092880188

# Row 12 = PORCELANATO
$wsCumplimiento.Range("D12").Value = 13501.66
$wsCumplimiento.Range("E12").Value = 48362.0603947566
$wsCumplimiento.Range("F12").Value = 0.2182484324228318

# Row 15 = TOTAL
$wsCumplimiento.Range("D15").Value = 26044.92
$wsCumplimiento.Range("E15").Value = 96009.91551083435
$wsCumplimiento.Range("F15").Value = 0.2133870394482494

# Column E width on sheet3 shrinks slightly (24 -> 23) as a side-effect of the content update.
# Note: the COM ColumnWidth setter applies pixel rounding/padding before it lands in the
# underlying column-width units, so we compensate (23 - 5/6) to land exactly on 23 once saved.
$wsCumplimiento.Columns("E").ColumnWidth = 22.16666666666667
